$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold plain text (e.g. "28.113.64", "0.4626", padded
# percentages like "  -2.79%  "). Force Text format first so Excel
# does not auto-convert these into numbers/percentages.

# --- Row 38 / 39 swap: FraxShare <-> VeChain ---
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02285"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.710"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.87%  "

# --- Price / Volume(1h) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.113.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4626"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3985"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08374"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.922.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.383"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.035"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06607"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.05%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.726"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.104.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.140.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.735"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9715"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09608"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.467"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.634"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.526"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.273"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06139"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6129"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1903"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.297"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5849"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.009"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.440"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06908"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.73"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.36%  "
